# Update "想去人数" (column F) values across the four sheets of the
# Guangzhou comic-con workbook, as produced by the latest data refresh.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (F column updates, by row)
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    2  = 829
    3  = 977
    4  = 767
    6  = 428
    7  = 663
    8  = 147
    9  = 1256
    10 = 683
    11 = 403
    13 = 178
    14 = 31
    15 = 829
    17 = 389
    18 = 368
    20 = 573
    21 = 131
    22 = 619
    24 = 908
    25 = 11
}
foreach ($row in $updates1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updates1[$row]
}

# Sheet "演出" (F column updates, by row)
$ws2 = $wb.Worksheets.Item("演出")
$updates2 = @{
    2  = 336
    5  = 639
    7  = 235
    11 = 106
}
foreach ($row in $updates2.Keys) {
    $ws2.Cells.Item($row, 6).Value = $updates2[$row]
}

# Sheet "本地生活" (F column updates, by row)
$ws3 = $wb.Worksheets.Item("本地生活")
$updates3 = @{
    2 = 373
}
foreach ($row in $updates3.Keys) {
    $ws3.Cells.Item($row, 6).Value = $updates3[$row]
}

# Sheet "全部类型" (F column updates, by row)
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    2  = 373
    3  = 336
    4  = 829
    5  = 977
    6  = 767
    8  = 428
    9  = 663
    10 = 147
    11 = 1256
    12 = 683
    15 = 403
    17 = 639
    18 = 178
    19 = 31
    20 = 829
    23 = 389
    24 = 368
    26 = 235
    28 = 573
    31 = 106
    32 = 106
    33 = 131
    34 = 619
    36 = 908
    37 = 11
}
foreach ($row in $updates4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updates4[$row]
}
